# The "展览" and "全部类型" sheets both list upcoming 漫展 (comic-con) events,
# one per row (row 1 is the header). The site's scrape that produced the
# "after" workbook dropped the oldest/cancelled event (old row 2,
# 2024-03-02 合肥·星芒1.5动漫嘉年华（取消）) and appended one new event at the
# bottom, while also refreshing a handful of "想去人数" (interest count)
# numbers on the events that stayed in the list. Net effect vs. the
# original rows 2-13:
#   - rows 3-13 slide up into rows 2-12 (B:I only - the row-number column A
#     keeps its original 1..11 values, only the trailing row 13 disappears)
#   - dimension shrinks from A1:I13 to A1:I12
#   - several F-column (想去人数) values get bumped to newer counts
#
# Both worksheets carry an identical copy of this table, so the same
# sequence is applied to each.

$wb = $excel.ActiveWorkbook

foreach ($sheetName in "展览", "全部类型") {
    $ws = $wb.Worksheets.Item($sheetName)

    # Slide the B:I content of rows 3-13 up into rows 2-12 (column A, the
    # row index, is left untouched so it keeps reading 1..11).
    $ws.Range("B3:I13").Copy()
    $ws.Range("B2").PasteSpecial()
    $excel.CutCopyMode = $false

    # The now-duplicated trailing row is removed, shrinking the used range
    # to A1:I12.
    $ws.Rows(13).Delete()

    # Refresh the "想去人数" counts that changed for the events which
    # remained on the list (row numbers below are the final, post-shift
    # positions).
    $ws.Range("F2").Value = 2738
    $ws.Range("F3").Value = 635
    $ws.Range("F5").Value = 6644
    $ws.Range("F6").Value = 1212
    $ws.Range("F8").Value = 20
    $ws.Range("F10").Value = 77
    $ws.Range("F11").Value = 10
}
